# Daily Report update: adds the 2026-02-12 (serial 46065) business day to
# Daily_Data, then refreshes the Today_Summary snapshot and the Monthly_Stats
# rollups so they reflect the newly appended day.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Daily_Data: append 24 rows (12 depositories x Registered/Eligible)
#    for the new date, starting right after the existing last row (97).
# ---------------------------------------------------------------------
$daily = $wb.Worksheets.Item("Daily_Data")

$newRows = @(
    @(46065, "ASAHI DEPOSITORY LLC Registered", 23953631.592, 0, 0, 0, 0, 23953631.592),
    @(46065, "ASAHI DEPOSITORY LLC Eligible", 2555897.608, 0, 0, 0, 0, 2555897.608),
    @(46065, "BRINK'S, INC. Registered", 16122359.646, 0, 0, 0, 0, 16122359.646),
    @(46065, "BRINK'S, INC. Eligible", 40640060.474, 0, 0, 0, 0, 40640060.474),
    @(46065, "CNT DEPOSITORY, INC. Registered", 12974598.079, 0, 0, 0, 0, 12974598.079),
    @(46065, "CNT DEPOSITORY, INC. Eligible", 15306765.903, 7944.85, 353108.525, -345163.675, 0, 14961602.228),
    @(46065, "DELAWARE DEPOSITORY Registered", 1552701.933, 0, 0, 0, 0, 1552701.933),
    @(46065, "DELAWARE DEPOSITORY Eligible", 16261227.756, 0, 11667.394, -11667.394, 0, 16249560.362),
    @(46065, "HSBC BANK, USA Registered", 3472271.68, 0, 0, 0, 0, 3472271.68),
    @(46065, "HSBC BANK, USA Eligible", 21150312.483, 0, 0, 0, 0, 21150312.483),
    @(46065, "INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Registered", 273789.87, 0, 0, 0, 0, 273789.87),
    @(46065, "INTERNATIONAL DEPOSITORY SERVICES OF DELAWARE Eligible", 3642206.244, 0, 0, 0, 0, 3642206.244),
    @(46065, "JP MORGAN CHASE BANK NA Registered", 12035257.32, 0, 0, 0, 0, 12035257.32),
    @(46065, "JP MORGAN CHASE BANK NA Eligible", 150151940.283, 0, 2441186.7, -2441186.7, 0, 147710753.583),
    @(46065, "LOOMIS INTERNATIONAL (US) LLC Registered", 7374299.767, 0, 0, 0, 0, 7374299.767),
    @(46065, "LOOMIS INTERNATIONAL (US) LLC Eligible", 23295383.436, 0, 0, 0, 0, 23295383.436),
    @(46065, "MALCA-AMIT ARMORED, INC. Registered", 0, 0, 0, 0, 0, 0),
    @(46065, "MALCA-AMIT ARMORED, INC. Eligible", 0, 0, 0, 0, 0, 0),
    @(46065, "MALCA-AMIT USA, LLC Registered", 1225506.264, 0, 0, 0, 0, 1225506.264),
    @(46065, "MALCA-AMIT USA, LLC Eligible", 798026.177, 0, 0, 0, 0, 798026.177),
    @(46065, "MANFRA, TORDELLA & BROOKES, LLC Registered", 6500477.621, 0, 0, 0, -130217.8, 6370259.821),
    @(46065, "MANFRA, TORDELLA & BROOKES, LLC Eligible", 12167803.719, 0, 0, 0, 130217.8, 12298021.519),
    @(46065, "STONEX PRECIOUS METALS LLC Registered", 7545291.14, 0, 0, 0, 0, 7545291.14),
    @(46065, "STONEX PRECIOUS METALS LLC Eligible", 233197.38, 0, 0, 0, 0, 233197.38)
)

$startRow = 98
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]

    $daily.Cells.Item($r, 1).Value = $row[0]
    $daily.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $daily.Cells.Item($r, 2).Value = $row[1]
    $daily.Cells.Item($r, 3).Value = $row[2]
    $daily.Cells.Item($r, 4).Value = $row[3]
    $daily.Cells.Item($r, 5).Value = $row[4]
    $daily.Cells.Item($r, 6).Value = $row[5]
    $daily.Cells.Item($r, 7).Value = $row[6]
    $daily.Cells.Item($r, 8).Value = $row[7]
}

# ---------------------------------------------------------------------
# 2) Today_Summary: a snapshot of the latest day's Eligible/Registered
#    totals per depository. Only the depositories whose latest-day
#    totals moved need updating (Total_Stock = Eligible + Registered).
# ---------------------------------------------------------------------
$today = $wb.Worksheets.Item("Today_Summary")

# Row 4: CNT DEPOSITORY, INC.
$today.Cells.Item(4, 2).Value = 14961602.228
$today.Cells.Item(4, 4).Value = 27936200.307

# Row 5: DELAWARE DEPOSITORY
$today.Cells.Item(5, 2).Value = 16249560.362
$today.Cells.Item(5, 4).Value = 17802262.295

# Row 8: JP MORGAN CHASE BANK NA
$today.Cells.Item(8, 2).Value = 147710753.583
$today.Cells.Item(8, 4).Value = 159746010.903

# Row 12: MANFRA, TORDELLA & BROOKES, LLC (both columns moved)
$today.Cells.Item(12, 2).Value = 12298021.519
$today.Cells.Item(12, 3).Value = 6370259.821
$today.Cells.Item(12, 4).Value = 18668281.34

# ---------------------------------------------------------------------
# 3) Monthly_Stats: roll the new day's RECEIVED/WITHDRAWN into the
#    month-to-date detail rows and refresh each TOTAL_TODAY (latest
#    day) column, then recompute the month's grand-total header row.
# ---------------------------------------------------------------------
$monthly = $wb.Worksheets.Item("Monthly_Stats")

# Row 11: CNT DEPOSITORY, INC. Eligible -- RECEIVED += 7944.85, WITHDRAWN += 353108.525
$monthly.Cells.Item(11, 3).Value = 7944.85
$monthly.Cells.Item(11, 4).Value = 3241301.818
$monthly.Cells.Item(11, 5).Value = 14961602.228

# Row 13: DELAWARE DEPOSITORY Eligible -- WITHDRAWN += 11667.394
$monthly.Cells.Item(13, 4).Value = 186666.419
$monthly.Cells.Item(13, 5).Value = 16249560.362

# Row 19: JP MORGAN CHASE BANK NA Eligible -- WITHDRAWN += 2441186.7
$monthly.Cells.Item(19, 4).Value = 9118134.300000001
$monthly.Cells.Item(19, 5).Value = 147710753.583

# Row 27: MANFRA, TORDELLA & BROOKES, LLC Eligible -- TOTAL_TODAY refresh
$monthly.Cells.Item(27, 4).Value = 849925.823
$monthly.Cells.Item(27, 5).Value = 12298021.519

# Row 28: MANFRA, TORDELLA & BROOKES, LLC Registered -- TOTAL_TODAY refresh
$monthly.Cells.Item(28, 5).Value = 6370259.821

# Row 2: month grand-total header (Eligible / Registered / Grand_Total)
$monthly.Cells.Item(2, 2).Value = 283535021.494
$monthly.Cells.Item(2, 3).Value = 92899967.112
$monthly.Cells.Item(2, 4).Value = 376434988.606
